$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rebalance the scoring: row 8's "Smart Score" changes from 1 to 0 ---
$ws.Range("D8").Value = 0

# --- Rename the table's old "Total" label (A19) to "Logic will be" first,
# so the shared string it occupies gets reused/renamed in place rather than
# leaving a stray "Total" entry behind. ---
$ws.Range("A19").Value = "Logic will be"

# --- Add the "Logic will be" explanation block in column G ---
$g2 = $ws.Range("G2")
$g2.Value = "Logic will be"
$g2.Font.Bold = $true
$g2.Font.Italic = $true

$ws.Range("G3").Value = "> 7 V + <7 S = V"
$ws.Range("G4").Value = "<7V + >7S = S"
$ws.Range("G5").Value = "<7V + <7S | >7V +>7S = Fail"
$ws.Range("G6").Value = "0V + 3S Accidental human easter egg"

# --- Remove the table's totals row (this drops the SUBTOTAL formulas that
# used to live in row 19 and shrinks the table back down to 18 data rows) ---
$lo = $ws.ListObjects.Item(1)
$lo.ShowTotals = $false
$ws.Range("A19:D19").ClearContents()

# --- Update the selected cell to match the new layout ---
$ws.Range("G16").Select()
